# Apply the diff described:
#  Summary sheet: B4 48->50, B5 2->0, B7 4->0
#  Symbols sheet: B2 48->50
#  Strategies sheet: D2 48->50, F2 2->0, G2 4->0.08333333333333333, H2 12->600,
#                     J2 text updated with new entry_ts/exit_ts (same day, 2-hour hold)

$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B4").Value = 50
$wsSummary.Range("B5").Value = 0
$wsSummary.Range("B7").Value = 0

# --- Symbols sheet ---
$wsSymbols = $wb.Worksheets.Item("Symbols")
$wsSymbols.Range("B2").Value = 50

# --- Strategies sheet ---
$wsStrategies = $wb.Worksheets.Item("Strategies")
$wsStrategies.Range("D2").Value = 50
$wsStrategies.Range("F2").Value = 0
$wsStrategies.Range("G2").Value = 0.08333333333333333
$wsStrategies.Range("H2").Value = 600
$wsStrategies.Range("J2").Value = "[{'strategy_name': 'Short Put', 'pnl': 50.0, 'entry_ts': '2025-01-01T10:00:00', 'exit_ts': '2025-01-01T12:00:00'}]"
